$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to a single blank space.
$ws.Name = " "

# Scroll the view back to the top-left corner (A1) instead of the
# previously-saved scroll position.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# "ficar sem cabecalho nem rodape" -> drop the printed footer text and
# shrink the now-unused bottom/footer page margins accordingly (the
# header itself, and its margin, are left untouched).
$ps = $ws.PageSetup
$ps.CenterFooter = ""
$ps.LeftFooter = ""
$ps.RightFooter = ""
$ps.BottomMargin = 56.699999999999996
$ps.FooterMargin = 36.850393700787386
